$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data keeps numeric-looking price strings (e.g. "1.00", "206.89")
# as literal text (trailing zeros, fixed decimals). Excel auto-converts such
# strings to numbers on assignment, so those specific cells are pre-formatted
# as Text ("@") to preserve the exact string, matching the source workbook.

# Row 2
$ws.Cells.Item(2, 4).Value = "81.422.31"
$ws.Cells.Item(2, 5).Value = "  +2.89%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.171.08"
$ws.Cells.Item(3, 5).Value = "  -0.49%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.17%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "206.89"
$ws.Cells.Item(5, 5).Value = "  -0.10%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "630.66"
$ws.Cells.Item(6, 5).Value = "  -0.11%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.291"
$ws.Cells.Item(7, 5).Value = "  +27.99%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "1.00"
$ws.Cells.Item(8, 5).Value = "  +0.05%  "

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.587"
$ws.Cells.Item(9, 5).Value = "  +1.93%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "3.172.11"
$ws.Cells.Item(10, 5).Value = "  -0.43%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.587"
$ws.Cells.Item(11, 5).Value = "  +1.47%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.0000259"
$ws.Cells.Item(12, 5).Value = "  +13.38%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.165"
$ws.Cells.Item(13, 5).Value = "  +1.69%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.34"
$ws.Cells.Item(14, 5).Value = "  -1.78%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.756.58"
$ws.Cells.Item(15, 5).Value = "  -0.25%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "31.80"
$ws.Cells.Item(16, 5).Value = "  +0.53%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "81.477.73"
$ws.Cells.Item(17, 5).Value = "  +3.28%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "3.174.66"
$ws.Cells.Item(18, 5).Value = "  -0.01%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "3.20"
$ws.Cells.Item(19, 5).Value = "  +12.96%  "

# Row 20
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "14.09"
$ws.Cells.Item(20, 5).Value = "  -2.11%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "9.19"
$ws.Cells.Item(21, 5).Value = "  -2.37%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "439.02"
$ws.Cells.Item(22, 5).Value = "  +1.39%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "5.15"
$ws.Cells.Item(23, 5).Value = "  +5.37%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "7.16"
$ws.Cells.Item(24, 5).Value = "  +4.97%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "5.32"
$ws.Cells.Item(25, 5).Value = "  +11.33%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "11.22"
$ws.Cells.Item(26, 5).Value = "  +1.47%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "3.347.70"
$ws.Cells.Item(27, 5).Value = "  +0.15%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "76.69"
$ws.Cells.Item(28, 5).Value = "  -0.02%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "1.00"
$ws.Cells.Item(29, 5).Value = "  -0.43%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "0.0000125"
$ws.Cells.Item(30, 5).Value = "  +8.12%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "9.13"
$ws.Cells.Item(31, 5).Value = "  +2.86%  "

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.999"
$ws.Cells.Item(32, 5).Value = "  -0.09%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "Bittensor"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "566.39"
$ws.Cells.Item(33, 5).Value = "  +9.38%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "Fetch.AI"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.51"
$ws.Cells.Item(34, 5).Value = "  +2.54%  "

# Row 35
$ws.Cells.Item(35, 2).Value = "Kaspa"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "0.153"
$ws.Cells.Item(35, 5).Value = "  +12.63%  "

# Row 36
$ws.Cells.Item(36, 2).Value = "PancakeSwap"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "2.03"
$ws.Cells.Item(36, 5).Value = "  +2.33%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +27.39%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "23.03"
$ws.Cells.Item(38, 5).Value = "  +1.88%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.00"
$ws.Cells.Item(39, 5).Value = "  +0.13%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.413"
$ws.Cells.Item(40, 5).Value = "  +3.87%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "RenderToken"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "6.00"
$ws.Cells.Item(41, 5).Value = "  +10.37%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "Stacks"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "2.04"
$ws.Cells.Item(42, 5).Value = "  +14.47%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "dogwifhat"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.07"
$ws.Cells.Item(43, 5).Value = "  +19.47%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +3.77%  "

# Row 45
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "159.64"
$ws.Cells.Item(45, 5).Value = "  -2.84%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  +0.01%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "189.40"
$ws.Cells.Item(47, 5).Value = "  -3.86%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  +3.31%  "

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "44.38"
$ws.Cells.Item(49, 5).Value = "  +3.27%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.782"
$ws.Cells.Item(50, 5).Value = "  -2.19%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "25.87"
$ws.Cells.Item(51, 5).Value = "  +6.28%  "
